$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 417, shifting existing rows 417:438 down to 418:439
$ws.Rows.Item(417).Insert()

# Populate the newly inserted row 417 with the new weekly data point
$ws.Range("A417").Value = 6
$ws.Range("B417").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C417").Value = "Metropolitana"
$ws.Range("D417").Value = 44706
$ws.Range("E417").Value = 13
$ws.Range("F417").Value = 100112039
$ws.Range("G417").Value = "Ciboulette"
$ws.Range("H417").Value = "Sin especificar"
$ws.Range("I417").Value = "Primera"
$ws.Range("J417").Value = 620
$ws.Range("K417").Value = 700
$ws.Range("L417").Value = 800
$ws.Range("M417").Value = 747
$ws.Range("N417").Value = "`$/docena de atados"
$ws.Range("O417").Value = "Región Metropolitana"
$ws.Range("P417").Value = 249
$ws.Range("Q417").Value = 3
$ws.Range("R417").Value = "Hortaliza"
